$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace "transmitter calibration" command with slope/offset cal variants,
# and add new "waypoints" command rows per updated serial ICD.
# (Order below matches the order new shared-string entries must be created in.)
$ws.Range("C17").Value = "waypoints"
$ws.Range("F17").Value = "read waypoints"
$ws.Range("C9").Value = "transmitter slope cal"
$ws.Range("C10").Value = "transmitter offset cal"
$ws.Range("F9").Value = "read transmitter slope values"
$ws.Range("F10").Value = "read transmitter offset values"
$ws.Range("C18").Value = "camera values"
$ws.Range("F18").Value = "read camera values"

# Update the active selection on the sheet to match the new edit location.
$ws.Range("F10").Select()
